# Adapt column header formatting to respective input file names:
# "<name>_old" -> "<name>_FV2210", "<name>_new" -> "<name>_FV2304"
# Then wrap the header/data range in a native Excel Table, and freeze the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the used range into a proper Excel Table ("Table1").
$range = $ws.Range("A1:U67")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"
$table.ShowTableStyleRowStripes = $true
$table.ShowTableStyleColumnStripes = $false

# Freeze the header row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
